# Weekly update: insert a new most-recent price record at row 17 for the
# "Fruta, Mercado Mayorista Lo Valledor de Santiago - Pomelo" sheet.
# This pushes all existing records (previously rows 17-49) down by one row
# (becoming rows 18-50) and populates the freshly opened row 17 with the
# new week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 17, shifting rows 17:49
# down to 18:50 (and growing the used range to A1:T50).
$ws.Rows.Item(17).Insert()

# Columns A, B, C, E, F, G, H, I, J, K are constant metadata for every
# record in this subset (market, region, product taxonomy, variety), so
# copy them straight from the row directly below (the record that used
# to be row 17 before the insert, now row 18).
for ($col = 1; $col -le 11; $col++) {
    $srcValue = $ws.Cells.Item(18, $col).Value2
    $ws.Cells.Item(17, $col).Value = $srcValue
}

# New week's own data.
$ws.Range("L17").Value = "Primera"
$ws.Range("D17").Value = 44935
$ws.Range("M17").Value = 210
$ws.Range("N17").Value = 10000
$ws.Range("O17").Value = 10000
$ws.Range("P17").Value = 10000
$ws.Range("Q17").Value = "$/caja 14 kilos"
$ws.Range("R17").Value = "Región de O'Higgins"
$ws.Range("S17").Value = 714
$ws.Range("T17").Value = 14
